# Implement json import and export
# Rename the sheet, add a "Number of Teacher" column before the existing
# "Number of students" column (renamed to "Number of Student"), and fill
# the new column with 0 for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab name) from "Schools" to "truong.xlsx"
$ws.Name = "truong.xlsx"

# 2. Insert a new column at C, shifting the old "Number of students" (C)
#    and "Address" (D) columns one place to the right (-> D, E).
$ws.Columns.Item(3).Insert()

# 3. Header row: new column C = "Number of Teacher"; rename old header
#    (now in D) from "Number of students" to "Number of Student".
$ws.Cells.Item(1, 3).Value = "Number of Teacher"
$ws.Cells.Item(1, 4).Value = "Number of Student"

# 4. Populate the new "Number of Teacher" column with 0 for every data row.
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 3).Value = 0

# 5. Column widths (best-fit) for the two "Number of ..." columns.
$ws.Columns.Item(3).ColumnWidth = 23
$ws.Columns.Item(4).ColumnWidth = 22.833333333333332
